$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two obsolete rows (old rows 20 and 21); this shifts the
# former row 22 up to become row 20, matching the new 19-data-row table.
$ws.Rows.Item(20).Delete()
$ws.Rows.Item(20).Delete()
